$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every data cell in this sheet stores plain text (prices/percent deltas are
# pre-formatted strings, not numbers/formulas). Some of the new values look like
# a plain decimal (e.g. "3.90"), which Excel would otherwise auto-convert to the
# number 3.9, silently dropping the trailing zero. Force text format for the
# assignment, then restore the default "Normal" style so no stray style index is
# left behind on cells that originally had none.
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "69.512.30"
$r.Style = "Normal"
$r = $ws.Range("E2")
$r.NumberFormat = "@"
$r.Value = "  +0.09%  "
$r.Style = "Normal"
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "2.494.34"
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = "  -0.71%  "
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = "  +0.01%  "
$r.Style = "Normal"
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "569.27"
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = "  -0.48%  "
$r.Style = "Normal"
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "166.16"
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = "  +0.24%  "
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = "  -0.03%  "
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = "  -0.23%  "
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = "  +0.56%  "
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = "  -2.69%  "
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = "  -0.12%  "
$r.Style = "Normal"
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "2.950.73"
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = "  -0.77%  "
$r.Style = "Normal"
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "69.363.78"
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = "  +0.05%  "
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = "  -0.08%  "
$r.Style = "Normal"
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "24.16"
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = "  -2.48%  "
$r.Style = "Normal"
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "2.432.64"
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = "  -3.34%  "
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = "  -0.51%  "
$r.Style = "Normal"
$r = $ws.Range("E19")
$r.NumberFormat = "@"
$r.Value = "  -1.23%  "
$r.Style = "Normal"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "352.64"
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = "  +1.14%  "
$r.Style = "Normal"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "3.90"
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = "  +0.11%  "
$r.Style = "Normal"
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "1.92"
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = "  -3.41%  "
$r.Style = "Normal"
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "69.38"
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = "  -1.26%  "
$r.Style = "Normal"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "3.80"
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = "  -2.59%  "
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = "  -1.30%  "
$r.Style = "Normal"
$r = $ws.Range("E27")
$r.NumberFormat = "@"
$r.Value = "  -2.22%  "
$r.Style = "Normal"
$r = $ws.Range("E28")
$r.NumberFormat = "@"
$r.Value = "  +1.09%  "
$r.Style = "Normal"
$r = $ws.Range("E29")
$r.NumberFormat = "@"
$r.Value = "  -1.56%  "
$r.Style = "Normal"
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "7.54"
$r.Style = "Normal"
$r = $ws.Range("E30")
$r.NumberFormat = "@"
$r.Value = "  -3.55%  "
$r.Style = "Normal"
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "3.52"
$r.Style = "Normal"
$r = $ws.Range("E31")
$r.NumberFormat = "@"
$r.Value = "  +139.12%  "
$r.Style = "Normal"
$r = $ws.Range("E32")
$r.NumberFormat = "@"
$r.Value = "  -3.12%  "
$r.Style = "Normal"
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "439.52"
$r.Style = "Normal"
$r = $ws.Range("E33")
$r.NumberFormat = "@"
$r.Value = "  -4.66%  "
$r.Style = "Normal"
$r = $ws.Range("E34")
$r.NumberFormat = "@"
$r.Value = "  +0.02%  "
$r.Style = "Normal"
$r = $ws.Range("E35")
$r.NumberFormat = "@"
$r.Value = "  -0.76%  "
$r.Style = "Normal"
$r = $ws.Range("E36")
$r.NumberFormat = "@"
$r.Value = "  -3.30%  "
$r.Style = "Normal"
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "153.13"
$r.Style = "Normal"
$r = $ws.Range("E37")
$r.NumberFormat = "@"
$r.Value = "  -2.63%  "
$r.Style = "Normal"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "19.07"
$r.Style = "Normal"
$r = $ws.Range("E38")
$r.NumberFormat = "@"
$r.Value = "  -0.07%  "
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = "  -1.59%  "
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = "  +0.03%  "
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = "  -1.04%  "
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = "  -2.15%  "
$r.Style = "Normal"
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = "  -1.91%  "
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = "  -2.34%  "
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = "  -3.65%  "
$r.Style = "Normal"
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "139.18"
$r.Style = "Normal"
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = "  -0.65%  "
$r.Style = "Normal"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.505"
$r.Style = "Normal"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = "  -2.46%  "
$r.Style = "Normal"
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = "  -0.94%  "
$r.Style = "Normal"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.574"
$r.Style = "Normal"
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = "  -0.77%  "
$r.Style = "Normal"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0925"
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = "  -0.39%  "
$r.Style = "Normal"
